# Commit: "remove links from workbooks to FixedIncome.xla"
#
# The workbook carries a single external reference (an [1]! link) to the
# FixedIncome.xla add-in, used by cell D8 ('General Settings') through the
# formula =[1]!qlSerializationPath(Trigger). We remove that link entirely
# and replace the formula result with a plain, hard-coded path string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General Settings")

$target  = $ws.Range("D8")
$scratch = $ws.Range("Z100")

# Write the new literal value into a throwaway cell that carries the
# default (unstyled) format, then copy only the VALUE back onto D8. This
# swaps the cell's content from the external-link formula to plain text
# without disturbing D8's own number format / style (otherwise a plain
# ".Value =" assignment on D8 would fork a brand-new style record).
$scratch.Value = "C:\Users\erik\junk\"
$scratch.Copy()
$target.PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false
$scratch.Clear()

# Drop the link to FixedIncome.xla (removes <externalReferences> from the
# workbook plus the xl/externalLinks/externalLink1.xml part).
$wb.BreakLink("FixedIncome.xla", 1)
